$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value. Using a loop + forcing
# the "Text" number format keeps these as strings (matching the
# source inline-string cells) instead of Excel auto-converting
# numeric-looking values (e.g. "8.44") or formatted strings
# (e.g. "45.140.63", "  +3.19%  ") into numbers/dates.
$updates = [ordered]@{
    'D2' = '45.140.63'
    'E2' = '  +3.19%  '
    'D3' = '2.366.60'
    'E3' = '  +1.18%  '
    'E4' = '  -0.13%  '
    'D5' = '310.64'
    'E5' = '  -0.42%  '
    'D6' = '107.69'
    'E6' = '  -1.43%  '
    'E7' = '  -0.16%  '
    'E8' = '  -0.08%  '
    'D9' = '0.615'
    'E9' = '  -0.67%  '
    'D10' = '40.79'
    'E10' = '  -0.84%  '
    'E11' = '  -0.49%  '
    'D12' = '8.44'
    'E12' = '  -2.02%  '
    'E13' = '  +1.16%  '
    'D14' = '0.974'
    'E14' = '  -3.37%  '
    'D15' = '2.725.07'
    'E15' = '  +1.33%  '
    'E16' = '  -1.85%  '
    'D17' = '2.377.39'
    'E17' = '  +1.52%  '
    'D18' = '45.127.70'
    'E18' = '  +3.34%  '
    'D19' = '14.59'
    'E19' = '  +10.95%  '
    'D20' = '7.26'
    'E20' = '  -4.59%  '
    'E21' = '  -1.14%  '
    'D22' = '73.10'
    'E22' = '  -1.57%  '
    'E23' = '  -0.62%  '
    'D24' = '259.15'
    'E24' = '  -3.88%  '
    'D25' = '2.28'
    'E25' = '  +0.19%  '
    'E26' = '  -0.25%  '
    'D27' = '11.08'
    'E27' = '  -0.91%  '
    'E28' = '  -5.94%  '
    'D29' = '2.34'
    'E29' = '  +2.34%  '
    'D30' = '0.0962'
    'E30' = '  +8.52%  '
    'D31' = '22.33'
    'E31' = '  -1.75%  '
    'D32' = '37.30'
    'E32' = '  -4.01%  '
    'D33' = '169.09'
    'E33' = '  +0.63%  '
    'E34' = '  +6.37%  '
    'E35' = '  -1.39%  '
    'E36' = '  +3.59%  '
    'B38' = 'NEARProtocol'
    'C38' = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
    'D38' = '3.90'
    'E38' = '  +1.79%  '
    'B39' = 'LidoDAOToken'
    'C39' = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
    'D39' = '2.93'
    'E39' = '  +2.07%  '
    'E40' = '  -3.77%  '
    'D41' = '1.74'
    'E41' = '  +0.96%  '
    'D42' = '99.61'
    'E42' = '  -4.80%  '
    'B43' = 'MultiversX'
    'C43' = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
    'D43' = '69.20'
    'E43' = '  -3.59%  '
    'B44' = 'Algorand'
    'C44' = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
    'D44' = '0.228'
    'E44' = '  -3.91%  '
    'D45' = '12.90'
    'E45' = '  -3.90%  '
    'E46' = '  -0.09%  '
    'D47' = '1.840.50'
    'E47' = '  +10.61%  '
    'D48' = '81.43'
    'E48' = '  +5.46%  '
    'D49' = '5.59'
    'E49' = '  +4.99%  '
    'D50' = '111.78'
    'E50' = '  -2.19%  '
    'D51' = '9.15'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
    $cell.Style = "Normal"
}
